$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the percentage profile row (row 3): B3:F3 = 0.2, 0.4, 0.6, 0.8, 1
$ws.Range("B3").Value = 0.2
$ws.Range("C3").Value = 0.4
$ws.Range("D3").Value = 0.6
$ws.Range("E3").Value = 0.8
$ws.Range("F3").Value = 1

# Add more "mgaye" file rows in column A: A9:A12 = 7, 8, 9, 10
$ws.Range("A9").Value = 7
$ws.Range("A10").Value = 8
$ws.Range("A11").Value = 9
$ws.Range("A12").Value = 10

# Update selection to match new active cell
$ws.Range("B12").Select()
